# B6-PowerPoint.pptx edit: swap the deck's colour theme palette for the
# "Office" palette (what used to live only in the notes-master theme part)
# and point the three data tables at the new table-style GUID.

$p = $ppt.ActivePresentation

# --- 1. Re-point the three tables (slides 14-16) at the new table style ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style.Name -eq "{E983B635-5CAB-4B85-8B89-86CB43C62EBA}") {
                $table.ApplyStyle("{C276418E-0801-4F4A-91C0-0546B1E0161D}", $false)
            }
        }
    }
}

# --- 2. Recolour the presentation theme to the Office palette -------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (PowerPoint RGB() packing,
# i.e. 0x00BBGGRR, matching the target a:srgbClr values below):
#   000000 FFFFFF 44546A E7E6E6 5B9BD5 ED7D31 A5A5A5 FFC000 4472C4 70AD47 0563C1 954F72
$officePalette = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officePalette[$i - 1]
}

Write-Host "theme recoloured and table styles updated"
